$d = $word.ActiveDocument

# 1) Locate the paragraph that currently holds "V32: " plus the _GoBack bookmark,
#    and replace its content (and the following, until-now-empty trailing
#    paragraph) with the new Vietnamese heading, a hyperlink placeholder,
#    four indented list items, and a paragraph that now carries the
#    relocated _GoBack bookmark.
$v32Index = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("V32:")) {
        $v32Index = $i
    }
}

$targetPara = $d.Paragraphs.Item($v32Index)
$r = $targetPara.Range
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:r><w:t xml:space="preserve">V32: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Cấp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>phát</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>và</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>giải</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>phóng</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bộ</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>nhớ</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>trong</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> C++</w:t></w:r></w:p><w:p></w:p><w:p><w:pPr><w:ind w:left="420"/></w:pPr><w:r><w:t xml:space="preserve">1/ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Cấp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>phát</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bộ</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>nhớ</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>động</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="420"/></w:pPr><w:r><w:t xml:space="preserve">2/ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Giải</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>phóng</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bộ</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>nhớ</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="420"/></w:pPr><w:r><w:t xml:space="preserve">3/ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Quan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>hệ</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>giữa</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> con </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>trỏ</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>với</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mảng</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="420"/></w:pPr><w:r><w:t xml:space="preserve">4/ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Thực</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>hành</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# 2) The paragraph right after the V32 heading is an empty placeholder for the
#    hyperlink. Use Hyperlinks.Add so Word mints the relationship and the
#    correct Hyperlink character style on the run.
$hpPara = $d.Paragraphs.Item($v32Index + 1)
$hlr = $hpPara.Range
$hlr.End = $hlr.End - 1
$d.Hyperlinks.Add($hlr, "https://youtu.be/mzfzAtXqkyM", [Type]::Missing, [Type]::Missing, "https://youtu.be/mzfzAtXqkyM")

# 3) Drop the old trailing empty paragraph that used to sit before the
#    sectPr — the relocated bookmark paragraph now plays that role.
$n = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($n)
if ($last.Range.Text -eq "") {
    $prev = $d.Paragraphs.Item($n - 1)
    $delRange = $d.Range($prev.Range.End - 1, $last.Range.End)
    $delRange.Delete()
}
